# Applies the "Updated cryptos list" data refresh to Sheet1 (cells B/C/D/E, rows 2-51).
# Numeric-looking Price values (column D) are written with a leading apostrophe so
# Excel keeps them as plain text (matching the original inlineStr/sharedString cells)
# instead of auto-coercing them into real numbers; the Style reset afterwards drops the
# "quote prefix" text-format style Excel applies, restoring the original (style-less) cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.519.46'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.914.93'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''245.27'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").Value = '''0.06733'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("D10").Value = '''110.94'
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = '''19.20'
$ws.Range("E11").Value = '  +3.99%  '
$ws.Range("D12").Value = '1.914.62'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '''0.07558'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = '''0.6684'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '''293.98'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '30.508.75'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '''12.96'
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '''0.000007589'
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D21").Value = '2.163.10'
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("D22").Value = '''5.491'
$ws.Range("E22").Value = '  +4.52%  '
$ws.Range("D23").Value = '''0.9998'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '''6.401'
$ws.Range("E24").Value = '  +3.04%  '
$ws.Range("D25").Value = '''9.476'
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("D26").Value = '''164.41'
$ws.Range("E26").Value = '  -2.53%  '
$ws.Range("D27").Value = '''20.47'
$ws.Range("E27").Value = '  -6.51%  '
$ws.Range("D28").Value = '''2.111'
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("D29").Value = '''0.1071'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +2.67%  '
$ws.Range("D31").Value = '''4.171'
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("E32").Value = '  +1.00%  '
$ws.Range("D33").Value = '''0.04977'
$ws.Range("D34").Value = '''0.7297'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("D36").Value = '''0.02052'
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '''2.732'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("D39").Value = '''2.673'
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '''110.96'
$ws.Range("E40").Value = '  +1.36%  '
$ws.Range("D41").Value = '''2.015'
$ws.Range("E41").Value = '  -2.30%  '
$ws.Range("D42").Value = '''0.4418'
$ws.Range("E42").Value = '  +3.65%  '
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = '''5.868'
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").Value = '''0.9995'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = '''68.36'
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("D47").Value = '''7.316'
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("D48").Value = '''49.10'
$ws.Range("E48").Value = '  -5.00%  '
$ws.Range("D49").Value = '''9.273'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '''0.1237'
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("D51").Value = '''0.2535'
$ws.Range("E51").Value = '  +3.67%  '

# Strip the auto-applied text NumberFormat/quote-prefix style so cells stay style-less,
# exactly like the rest of the sheet.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
